$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Column D holds price strings formatted with thousands separators (e.g. "64.504.38")
# that Excel's COM Value setter would otherwise auto-convert to numbers. Force the
# whole data range to Text format first so the literal strings are preserved verbatim.
$ws.Range("D2:D51").NumberFormat = "@"

$ws.Range("D2").Value = "64.504.38"
$ws.Range("E2").Value = "  -2.61%  "

$ws.Range("D3").Value = "3.178.68"
$ws.Range("E3").Value = "  -4.09%  "

$ws.Range("E4").Value = "  +0.01%  "

$ws.Range("D5").Value = "570.79"
$ws.Range("E5").Value = "  -2.79%  "

$ws.Range("D6").Value = "169.53"
$ws.Range("E6").Value = "  -7.59%  "

$ws.Range("D7").Value = "0.609"
$ws.Range("E7").Value = "  -6.20%  "

$ws.Range("D8").Value = "0.999"
$ws.Range("E8").Value = "  -0.16%  "

$ws.Range("D9").Value = "3.186.90"
$ws.Range("E9").Value = "  -3.83%  "

$ws.Range("E10").Value = "  -3.58%  "

$ws.Range("D11").Value = "6.82"
$ws.Range("E11").Value = "  +0.20%  "

$ws.Range("D12").Value = "0.387"
$ws.Range("E12").Value = "  -3.44%  "

$ws.Range("D13").Value = "3.738.92"
$ws.Range("E13").Value = "  -3.96%  "

$ws.Range("E14").Value = "  -2.25%  "

$ws.Range("D15").Value = "64.524.57"
$ws.Range("E15").Value = "  -2.65%  "

$ws.Range("D16").Value = "25.47"
$ws.Range("E16").Value = "  -2.50%  "

$ws.Range("E17").Value = "  -2.41%  "

$ws.Range("D18").Value = "3.172.80"
$ws.Range("E18").Value = "  -5.61%  "

$ws.Range("D19").Value = "419.99"
$ws.Range("E19").Value = "  -1.36%  "

$ws.Range("D20").Value = "12.94"
$ws.Range("E20").Value = "  -1.81%  "

$ws.Range("D21").Value = "5.37"
$ws.Range("E21").Value = "  -2.77%  "

$ws.Range("E22").Value = "  -3.55%  "

$ws.Range("D23").Value = "1.00"
$ws.Range("E23").Value = "  -0.04%  "

$ws.Range("E24").Value = "  -0.10%  "

$ws.Range("D25").Value = "70.40"
$ws.Range("E25").Value = "  -2.07%  "

$ws.Range("E26").Value = "  +0.48%  "

$ws.Range("D27").Value = "0.490"
$ws.Range("E27").Value = "  -4.78%  "

$ws.Range("E28").Value = "  -6.44%  "

$ws.Range("D29").Value = "8.82"
$ws.Range("E29").Value = "  -1.35%  "

$ws.Range("E30").Value = "  +0.13%  "

$ws.Range("E31").Value = "  -4.87%  "

$ws.Range("D32").Value = "21.80"
$ws.Range("E32").Value = "  -2.41%  "

$ws.Range("E33").Value = "  -0.10%  "

$ws.Range("D34").Value = "5.10"
$ws.Range("E34").Value = "  -1.47%  "

$ws.Range("D35").Value = "6.34"
$ws.Range("E35").Value = "  -3.56%  "

$ws.Range("E36").Value = "  -3.31%  "

$ws.Range("D37").Value = "157.64"
$ws.Range("E37").Value = "  -1.39%  "

$ws.Range("E38").Value = "  -4.80%  "

$ws.Range("D39").Value = "2.732.00"
$ws.Range("E39").Value = "  -5.13%  "

$ws.Range("E40").Value = "  -4.94%  "

$ws.Range("D41").Value = "24.44"
$ws.Range("E41").Value = "  -7.40%  "

$ws.Range("D42").Value = "4.20"
$ws.Range("E42").Value = "  -2.65%  "

$ws.Range("D43").Value = "39.20"
$ws.Range("E43").Value = "  -2.24%  "

$ws.Range("D44").Value = "0.714"
$ws.Range("E44").Value = "  -6.49%  "

$ws.Range("D45").Value = "0.0625"
$ws.Range("E45").Value = "  -5.76%  "

$ws.Range("D46").Value = "5.72"
$ws.Range("E46").Value = "  -3.15%  "

$ws.Range("D47").Value = "0.0265"
$ws.Range("E47").Value = "  -2.57%  "

# Row 48 now holds InjectiveProtocol (previously row 49's data)
$ws.Range("B48").Value = "InjectiveProtocol"
$ws.Range("C48").Value = "https://coinranking.com/coin/PkY9BmsyW+injectiveprotocol-inj"
$ws.Range("D48").Value = "21.73"
$ws.Range("E48").Value = "  -6.53%  "

# Row 49 now holds Bittensor (previously row 48's data)
$ws.Range("B49").Value = "Bittensor"
$ws.Range("C49").Value = "https://coinranking.com/coin/pgv7xSFi6+bittensor-tao"
$ws.Range("D49").Value = "294.65"
$ws.Range("E49").Value = "  -5.91%  "

# Row 50 now holds dogwifhat (previously row 51's data)
$ws.Range("B50").Value = "dogwifhat"
$ws.Range("C50").Value = "https://coinranking.com/coin/sZUrmToWF+dogwifhat-wif"
$ws.Range("D50").Value = "2.01"
$ws.Range("E50").Value = "  -12.67%  "

# Row 51 now holds FirstDigitalUSD (previously row 50's data)
$ws.Range("B51").Value = "FirstDigitalUSD"
$ws.Range("C51").Value = "https://coinranking.com/coin/cpjRxjFYD+firstdigitalusd-fdusd"
$ws.Range("D51").Value = "0.997"
$ws.Range("E51").Value = "  -0.28%  "
